# Generate Report for Handoff
# Update localization status from "In Translation" to "Ready for handoff"
# and bump the associated timestamps forward by 30 seconds, on all three
# worksheets (Overview, zh-cn, de-de). Also widen the "Status" columns to
# fit the new, longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 19:02:57"

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 19:02:52"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 19:02:57"

# Widen the Status columns to accommodate the longer "Ready for handoff" text
$wsOverview.Columns.Item(5).ColumnWidth = 16.3826548258464
$wsOverview.Columns.Item(6).ColumnWidth = 16.3826548258464
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3826548258464
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3826548258464
